# Trade #35 closed at 2026-02-17 20:57:59 - unknown UNKNOWN +0.000%
#
# Updates:
#  - Summary sheet: refreshed aggregate metrics (capital, P&L, trade counts, win rate)
#  - Strategy Status sheet: refreshed MarketMaking strategy row
#  - All Trades sheet: closed trade #63 (row 64) and appended new open trade #96 (row 97)
#  - MarketMaking sheet: closed trade #63 (row 31) and appended new open trade #96 (row 64)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1400.41
$summary.Range("B4").Value = 0.21
$summary.Range("B5").Value = 0.07000000000000001
$summary.Range("B6").Value = 63
$summary.Range("B8").Value = 26
$summary.Range("B9").Value = 46.03

# ---------------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 100.41
$status.Range("D5").Value = 30
$status.Range("E5").Value = 0.1
$status.Range("F5").Value = 0.41
$status.Range("G5").Value = 50

# ---------------------------------------------------------------------------
# All Trades sheet
# Columns: A Trade#, B Date, C Time, D Strategy, E Side, F Entry Price,
#          G Exit Price, H Status, I P&L %, J P&L $, K Capital After,
#          L Exit Reason, M Duration (min), N Entry Slippage (bps),
#          O Exit Slippage (bps), P Confidence, Q Entry Reason
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Close out existing trade #63 (row 64)
$allTrades.Cells.Item(64, 7).Value = 0.045539
$allTrades.Cells.Item(64, 8).Value = "CLOSED"
$allTrades.Cells.Item(64, 9).Value = -43.0763
$allTrades.Cells.Item(64, 10).Value = -0.03
$allTrades.Cells.Item(64, 11).Value = 100.41
$allTrades.Cells.Item(64, 12).Value = "early_exit"
$allTrades.Cells.Item(64, 13).Value = 0.14

# Append new trade #96 (row 97)
$allTrades.Cells.Item(97, 1).Value = 96
$allTrades.Cells.Item(97, 2).NumberFormat = "@"
$allTrades.Cells.Item(97, 2).Value = "2026-02-17"
$allTrades.Cells.Item(97, 3).NumberFormat = "@"
$allTrades.Cells.Item(97, 3).Value = "20:57:53"
$allTrades.Cells.Item(97, 4).Value = "MarketMaking"
$allTrades.Cells.Item(97, 5).Value = "UP"
$allTrades.Cells.Item(97, 6).Value = 0.08
$allTrades.Cells.Item(97, 7).Value = ""
$allTrades.Cells.Item(97, 8).Value = "OPEN"
$allTrades.Cells.Item(97, 9).Value = 0
$allTrades.Cells.Item(97, 10).Value = 0
$allTrades.Cells.Item(97, 11).Value = 100.4455022889912
$allTrades.Cells.Item(97, 12).Value = ""
$allTrades.Cells.Item(97, 13).Value = 0
$allTrades.Cells.Item(97, 14).Value = 0
$allTrades.Cells.Item(97, 15).Value = 0
$allTrades.Cells.Item(97, 16).Value = 0.6
$allTrades.Cells.Item(97, 17).Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------------
# MarketMaking sheet
# Columns: A Trade#, B Date, C Time, D Strategy, E Side, F Entry Price,
#          G Exit Price, H Status, I P&L %, J P&L $, K Capital After,
#          L Entry Slippage (bps), M Exit Slippage (bps), N Confidence,
#          O Entry Reason, P Exit Reason, Q Duration (min)
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")

# Close out existing trade #63 (row 31)
$marketMaking.Cells.Item(31, 7).Value = 0.045539
$marketMaking.Cells.Item(31, 8).Value = "CLOSED"
$marketMaking.Cells.Item(31, 9).Value = -43.0763
$marketMaking.Cells.Item(31, 10).Value = -0.03
$marketMaking.Cells.Item(31, 11).Value = 100.41
$marketMaking.Cells.Item(31, 16).Value = "early_exit"
$marketMaking.Cells.Item(31, 17).Value = 0.14

# Append new trade #96 (row 64)
$marketMaking.Cells.Item(64, 1).Value = 96
$marketMaking.Cells.Item(64, 2).NumberFormat = "@"
$marketMaking.Cells.Item(64, 2).Value = "2026-02-17"
$marketMaking.Cells.Item(64, 3).NumberFormat = "@"
$marketMaking.Cells.Item(64, 3).Value = "20:57:53"
$marketMaking.Cells.Item(64, 4).Value = "MarketMaking"
$marketMaking.Cells.Item(64, 5).Value = "UP"
$marketMaking.Cells.Item(64, 6).Value = 0.08
$marketMaking.Cells.Item(64, 7).Value = ""
$marketMaking.Cells.Item(64, 8).Value = "OPEN"
$marketMaking.Cells.Item(64, 9).Value = 0
$marketMaking.Cells.Item(64, 10).Value = 0
$marketMaking.Cells.Item(64, 11).Value = 100.4455022889912
$marketMaking.Cells.Item(64, 12).Value = 0
$marketMaking.Cells.Item(64, 13).Value = 0
$marketMaking.Cells.Item(64, 14).Value = 0.6
$marketMaking.Cells.Item(64, 15).Value = "Normal spread capture: 19600 bps"
$marketMaking.Cells.Item(64, 16).Value = ""
$marketMaking.Cells.Item(64, 17).Value = 0
